# "Admin" sheet: rebuild the task/report layout to make room for the new
# API Errors / Succesful Requests / Server Exceptions / Bets rows and the
# second "Book2" report block, then add the new "Admin TODO" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the old rows so the sheet can be rebuilt with the new structure
# (rows shift around quite a bit, so start from a clean slate).
$ws.Rows("1:36").Delete()

# New sheet for the admin TODO list, placed right after "Admin".
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Admin TODO"

# --- Admin Modules header ---
$ws.Range("A1").Value = "Admin Modules"
$ws.Range("A2").Value = "Dashboard:"
$ws.Range("E2").Value = "Type"
$ws.Range("B3").Value = "API Stats"
$ws.Range("E3").Value = "Table"

# --- Admin TODO sheet: Todo list ---
$ws2.Range("A1").Value = "Todo:"

# --- API Stats table rows (new metrics inserted above the existing ones) ---
$ws.Range("C7").Value = "Bets"
$ws.Range("C5").Value = "API Errors"
$ws.Range("C4").Value = "API Succesful Requests"
$ws.Range("C6").Value = "Server Exceptions (HTTP 500)"
$ws.Range("C8").Value = "Suggestions"
$ws.Range("C9").Value = "Admins"
$ws.Range("C10").Value = "Users"
$ws.Range("C11").Value = "Books"

$ws2.Range("A2").Value = "Dashboard"

# --- Book Stats block ---
$ws.Range("B13").Value = "Book Stats"
$ws.Range("C14").Value = "Book1"
$ws.Range("D15").Value = "API Succesful Requests"
$ws.Range("D16").Value = "API Errors"
$ws.Range("D17").Value = "Bets"
$ws.Range("D18").Value = "Suggestions"
$ws.Range("C19").Value = "Book2"
$ws.Range("D20").Value = "API Succesful Requests"
$ws.Range("D21").Value = "API Errors"
$ws.Range("D22").Value = "Bets"
$ws.Range("D23").Value = "Suggestions"

# --- API Activity charts ---
$ws.Range("B25").Value = "API Activity Chart last 24 hrs"
$ws.Range("E25").Value = "line area chart"
$ws.Range("B26").Value = "API Activity Chart last 7 days"
$ws.Range("E26").Value = "line area chart"
$ws.Range("B27").Value = "API Activity Chart this month"
$ws.Range("E27").Value = "line area chart"

# --- Admins / Users / Books / Reports sections ---
$ws.Range("A29").Value = "Admins:"
$ws.Range("B30").Value = "Page to add or remove admins"

$ws.Range("A32").Value = "Users:"
$ws.Range("B33").Value = "Page to add or remove users"
$ws.Range("B34").Value = "Page to add a book to a certain user"

$ws.Range("A36").Value = "Books:"
$ws.Range("B37").Value = "Page to add, edit/config or disable books"

$ws.Range("A39").Value = "Reports:"
$ws.Range("B40").Value = "API Stats"
$ws.Range("B41").Value = "Book Stats"
$ws.Range("B42").Value = "Financials"
$ws.Range("B43").Value = "Suggestions Tracker"

# --- Admin TODO sheet: rest of the list ---
$ws2.Range("A3").Value = "Admins"
$ws2.Range("A4").Value = "Users"
$ws2.Range("A5").Value = "Books"
$ws2.Range("A6").Value = "Reports"

# Column C on "Admin" needs to be widened for the longer labels now in it.
$ws.Columns.Item(3).ColumnWidth = 21.6

# Column A on "Admin TODO" sized to fit its labels.
$ws2.Columns.Item(1).ColumnWidth = 10.14

# Page setup (portrait) was (re)applied to "Admin" when it was reprinted.
$ws.PageSetup.Orientation = 1

# Selection / active-tab state from the authored workbook.
$ws.Range("D15").Select() | Out-Null
$ws2.Range("B1").Select() | Out-Null
